function Set-HitsRow {
    param($ws, $row, $hits, $pct)
    $ws.Range("B$row").Value = $hits
    $dCell = $ws.Range("D$row")
    $dCell.NumberFormat = "@"
    $dCell.Value = $pct
}

$wb = $excel.ActiveWorkbook

# --- Sheet: Total Hits ---
$ws = $wb.Worksheets.Item("Total Hits")
Set-HitsRow $ws 2 1839 "63.15%"
Set-HitsRow $ws 3 3701 "63.55%"
Set-HitsRow $ws 4 5572 "63.78%"
Set-HitsRow $ws 5 7449 "63.95%"
Set-HitsRow $ws 6 9303 "63.89%"

# --- Sheet: Hits_entity ---
$ws = $wb.Worksheets.Item("Hits_entity")
Set-HitsRow $ws 2 912  "64.64%"
Set-HitsRow $ws 3 1858 "65.84%"
Set-HitsRow $ws 4 2792 "65.96%"
Set-HitsRow $ws 5 3727 "66.03%"
Set-HitsRow $ws 6 4653 "65.95%"

# --- Sheet: Hits_numerical ---
$ws = $wb.Worksheets.Item("Hits_numerical")
Set-HitsRow $ws 2 338  "51.60%"
Set-HitsRow $ws 3 678  "51.76%"
Set-HitsRow $ws 5 1393 "53.17%"

# --- Sheet: Hits_boolean ---
$ws = $wb.Worksheets.Item("Hits_boolean")
Set-HitsRow $ws 2 399  "69.63%"
Set-HitsRow $ws 3 787  "68.67%"
Set-HitsRow $ws 4 1174 "68.30%"
Set-HitsRow $ws 5 1570 "68.50%"
Set-HitsRow $ws 6 1966 "68.62%"
